# pipe heat loss calculation updated; views added and updated; Code reformatted
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet/tab (Tabelle1 -> Parameter) ---
$ws.Name = "Parameter"

# --- Drop the now-unused last row (table shrinks from 50 to 49 rows) ---
[void]$ws.Rows.Item(50).Delete()

# --- Rewrite the "Netzauslegung" network-design block (rows 34-40): ---
# split the old single Vor-/Ruecklauftemperatur pair into separate
# winter/summer pairs, add ground temperatures winter/summer, and move
# "Gleichzeitigkeit" down into this block.
$ws.Range("A34").Value = "Netzauslegung"
$ws.Range("B34").Value = "Vorlauftemperatur Winter (K)"
$ws.Range("C34").Value = 348.15
$ws.Range("C34").NumberFormat = "0.00"
$ws.Range("D34").Value = "75°C Vorlauftemperatur Auslegung"

$ws.Range("A35").Value = "Netzauslegung"
$ws.Range("B35").Value = "Rücklauftemperatur Winter (K)"
$ws.Range("C35").Value = 328.15
$ws.Range("C35").NumberFormat = "0.00"
$ws.Range("D35").Value = "55°C Rücklauftemperatur Auslegung"

$ws.Range("A36").Value = "Netzauslegung"
$ws.Range("B36").Value = "Vorlauftemperatur Sommer (K)"
$ws.Range("C36").Value = 338.15
$ws.Range("C36").NumberFormat = "0.00"
$ws.Range("D36").Value = "65°C Vorlauftemperatur Sommer"

$ws.Range("A37").Value = "Netzauslegung"
$ws.Range("B37").Value = "Rücklauftemperatur Sommer (K)"
$ws.Range("C37").Value = 328.15
$ws.Range("C37").NumberFormat = "0.00"
$ws.Range("D37").Value = "55°C Rücklauftemperatur Sommer"

$ws.Range("A38").Value = "Netzauslegung"
$ws.Range("B38").Value = "Bodentemperatur Winter (K)"
$ws.Range("C38").Value = 276.64999999999998
$ws.Range("C38").NumberFormat = "0.00"
$ws.Range("D38").Value = "3,5°C Bodentemperaturen "

$ws.Range("A39").Value = "Netzauslegung"
$ws.Range("B39").Value = "Bodentemperatur Sommer (K)"
$ws.Range("C39").Value = 290.64999999999998
$ws.Range("C39").NumberFormat = "0.00"
$ws.Range("D39").ClearContents()
$ws.Range("E39").ClearContents()

$ws.Range("A40").Value = "Netzauslegung"
$ws.Range("B40").Value = "Gleichzeitigkeit (80%)"
$ws.Range("C40").Value = 0.8
$ws.Range("C40").NumberFormat = "0.0"
$ws.Range("D40").Value = "Vorgabe der Gleichzeitigkeit für die Berechnung des Leistungsbedarfs"
$ws.Range("E40").Value = "Vorlesung"

# --- "Netzauslegung_initial" block (rows 41-45): unchanged data, shifted down one row ---
$ws.Range("A41").Value = "Netzauslegung_initial"
$ws.Range("B41").Value = "Strömungsgeschwindigkeit (m/s)"
$ws.Range("C41").Value = 3
$ws.Range("C41").NumberFormat = "0.0"

$ws.Range("A42").Value = "Netzauslegung_initial"
$ws.Range("B42").Value = "Dichte Wasser Vorlauf (kg/m³)"
$ws.Range("C42").Value = 974.98900000000003
$ws.Range("C42").NumberFormat = "0.0000"
$ws.Range("D42").Value = "4 bar, 75°C peacesoftware.de Auf Druck und Temperatur anpassen!"

$ws.Range("A43").Value = "Netzauslegung_initial"
$ws.Range("B43").Value = "Dichte Wasser Rücklauf (kg/m³)"
$ws.Range("C43").Value = 985.83690000000001
$ws.Range("C43").NumberFormat = "0.0000"
$ws.Range("D43").Value = "4 bar, 55°C peacesoftware.de Auf Druck und Temperatur anpassen!"

$ws.Range("A44").Value = "Netzauslegung_initial"
$ws.Range("B44").Value = "Kinematische Viskosität Vorlauf (m²/s)"
$ws.Range("C44").Formula = "=0.38751610493371*0.00001"
$ws.Range("C44").NumberFormat = "0.0000000"

$ws.Range("A45").Value = "Netzauslegung_initial"
$ws.Range("B45").Value = "Kinematische Viskosität Rücklauf (m²/s)"
$ws.Range("C45").Formula = "=0.51128491182691*0.00001"
$ws.Range("C45").NumberFormat = "0.0000000"

# --- "#Rohrparameter" section header (row 46) ---
$ws.Range("A46").Value = "#Rohrparameter"
$ws.Range("B46").ClearContents()
$ws.Range("C46").ClearContents()

# --- "Rohrauswahl" block (rows 47-49) ---
$ws.Range("A47").Value = "Rohrauswahl"
$ws.Range("B47").Value = "Initiale Dämmung"
$ws.Range("C47").Value = "DS1"

$ws.Range("A48").Value = "Rohrauswahl"
$ws.Range("B48").Value = "Rohrrauheit k (mm)"
$ws.Range("C48").Value = 0.01
$ws.Range("C48").NumberFormat = "0.00"

$ws.Range("A49").Value = "Rohrauswahl"
$ws.Range("B49").Value = "Untergrenze Hydraulisch glatt"
$ws.Range("C49").Value = 2320
$ws.Range("C49").NumberFormat = "0"

# --- Resize the "Tabelle1" table object to the new, smaller range ---
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("A1:E49"))

# --- Update the view: select row 47 (as last edited by the author) ---
[void]$ws.Rows.Item(47).Select()
